# OW-248 externalized the graph data files into acuo-data
#
# 1) The placeholder Portfolio ID "acc1" used throughout every sheet's
#    column B is renamed to the real externalized account id "ACUOSG8745".
# 2) The workbook's active/selected tab moves from "OIS-Cleared" back to
#    "IRS-Cleared", and the OIS-Cleared sheet's selection/scroll resets to
#    the top of the sheet (A1 / B2) since it is no longer the focused tab.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the placeholder account id everywhere it appears ---------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("acc1", "ACUOSG8745")
}

# --- 2. Reset view state on the sheets that had their scroll/selection --
#        repositioned in the authored edit.

# OIS-Cleared loses focus; its window scrolls back to the top-left and the
# selection moves off the old AA2 cell to B2.
$wsOis = $wb.Worksheets.Item("OIS-Cleared")
$wsOis.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsOis.Range("B2").Select()

# IRS-Bilateral keeps its own selection (AE90) but its scroll position
# moves from R82 to A52.
$wsBilateral = $wb.Worksheets.Item("IRS-Bilateral")
$wsBilateral.Activate()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
$wsBilateral.Range("AE90").Select()

# IRS-Cleared becomes the active tab again; its selection stays on P137
# but the scroll position moves from Z286 to A99.
$wsIrs = $wb.Worksheets.Item("IRS-Cleared")
$wsIrs.Activate()
$excel.ActiveWindow.ScrollRow = 99
$excel.ActiveWindow.ScrollColumn = 1
$wsIrs.Range("P137").Select()
